$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Tipo" header in D1, matching the style used by the other headers
$ws.Range("D1").Value = "Tipo"
$ws.Range("A1").Copy()
$ws.Range("D1").PasteSpecial(-4122)

# Update MSE (column B) and R2 (column C) values
$ws.Range("B2").Value = 0.0531332828077846
$ws.Range("C2").Value = 0.9984374408005952

$ws.Range("B3").Value = 0.06076159267096478
$ws.Range("C3").Value = 0.9994236490633295

$ws.Range("B4").Value = 0.05595975722587399
$ws.Range("C4").Value = 0.9992423394300239

# Add new "Tipo" column values
$ws.Range("D2").Value = "multiple"
$ws.Range("D3").Value = "multiple"
$ws.Range("D4").Value = "multiple"
